$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.548.00'
$ws.Range('E2').Value = '  +0.87%  '
$ws.Range('D3').Value = '1.597.53'
$ws.Range('E3').Value = '  -1.14%  '
$ws.Range('E4').Value = '  +0.92%  '
$origStyle_D5 = $ws.Range('D5').Style
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '208.44'
$ws.Range('D5').Style = $origStyle_D5
$ws.Range('E5').Value = '  -0.23%  '
$ws.Range('E6').Value = '  -3.96%  '
$ws.Range('E7').Value = '  +0.77%  '
$origStyle_D8 = $ws.Range('D8').Style
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '22.31'
$ws.Range('D8').Style = $origStyle_D8
$ws.Range('E8').Value = '  -3.67%  '
$ws.Range('E9').Value = '  -1.23%  '
$origStyle_D10 = $ws.Range('D10').Style
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0589'
$ws.Range('D10').Style = $origStyle_D10
$ws.Range('E10').Value = '  -3.06%  '
$origStyle_D11 = $ws.Range('D11').Style
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0870'
$ws.Range('D11').Style = $origStyle_D11
$ws.Range('E11').Value = '  -0.54%  '
$ws.Range('D12').Value = '1.823.96'
$ws.Range('D13').Value = '1.592.24'
$ws.Range('E13').Value = '  -1.81%  '
$ws.Range('E14').Value = '  -3.57%  '
$origStyle_D15 = $ws.Range('D15').Style
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.540'
$ws.Range('D15').Style = $origStyle_D15
$ws.Range('E15').Value = '  -2.98%  '
$origStyle_D16 = $ws.Range('D16').Style
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '63.55'
$ws.Range('D16').Style = $origStyle_D16
$ws.Range('E16').Value = '  -2.22%  '
$ws.Range('D17').Value = '27.543.77'
$ws.Range('E17').Value = '  -0.13%  '
$origStyle_D18 = $ws.Range('D18').Style
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '216.64'
$ws.Range('D18').Style = $origStyle_D18
$ws.Range('E18').Value = '  -4.81%  '
$origStyle_D19 = $ws.Range('D19').Style
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '7.41'
$ws.Range('D19').Style = $origStyle_D19
$ws.Range('E19').Value = '  -1.95%  '
$ws.Range('D20').Value = '0.0₃0689'
$ws.Range('E20').Value = '  -3.94%  '
$ws.Range('E21').Value = '  +0.47%  '
$ws.Range('E22').Value = '  -1.91%  '
$origStyle_D23 = $ws.Range('D23').Style
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '9.77'
$ws.Range('D23').Style = $origStyle_D23
$ws.Range('E23').Value = '  -2.44%  '
$ws.Range('E24').Value = '  -0.87%  '
$origStyle_D25 = $ws.Range('D25').Style
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '154.25'
$ws.Range('D25').Style = $origStyle_D25
$ws.Range('E25').Value = '  +0.34%  '
$ws.Range('E26').Value = '  +0.63%  '
$ws.Range('E27').Value = '  -2.52%  '
$ws.Range('E28').Value = '  -2.62%  '
$ws.Range('E29').Value = '  -4.38%  '
$ws.Range('E30').Value = '  -1.42%  '
$origStyle_D31 = $ws.Range('D31').Style
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.0466'
$ws.Range('D31').Style = $origStyle_D31
$ws.Range('E31').Value = '  -2.61%  '
$ws.Range('E32').Value = '  -2.77%  '
$ws.Range('D33').Value = '1.370.09'
$ws.Range('E33').Value = '  -0.82%  '
$ws.Range('E34').Value = '  -3.15%  '
$ws.Range('E35').Value = '  -0.64%  '
$ws.Range('B36').Value = 'TrustWalletToken'
$ws.Range('C36').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$origStyle_D36 = $ws.Range('D36').Style
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.960'
$ws.Range('D36').Style = $origStyle_D36
$ws.Range('E36').Value = '  -3.40%  '
$ws.Range('B37').Value = 'HuobiToken'
$ws.Range('C37').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$origStyle_D37 = $ws.Range('D37').Style
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.32'
$ws.Range('D37').Style = $origStyle_D37
$ws.Range('E37').Value = '  -0.52%  '
$ws.Range('E38').Value = '  -2.23%  '
$origStyle_D39 = $ws.Range('D39').Style
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.535'
$ws.Range('D39').Style = $origStyle_D39
$ws.Range('E39').Value = '  -3.03%  '
$origStyle_D40 = $ws.Range('D40').Style
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.815'
$ws.Range('D40').Style = $origStyle_D40
$ws.Range('E40').Value = '  -3.82%  '
$ws.Range('E41').Value = '  +0.61%  '
$origStyle_D42 = $ws.Range('D42').Style
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.968'
$ws.Range('D42').Style = $origStyle_D42
$ws.Range('E42').Value = '  -5.37%  '
$origStyle_D43 = $ws.Range('D43').Style
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '5.34'
$ws.Range('D43').Style = $origStyle_D43
$ws.Range('E43').Value = '  -2.03%  '
$origStyle_D44 = $ws.Range('D44').Style
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '63.93'
$ws.Range('D44').Style = $origStyle_D44
$ws.Range('E44').Value = '  -2.09%  '
$ws.Range('E45').Value = '  -2.96%  '
$ws.Range('D46').Value = '1.734.77'
$ws.Range('E46').Value = '  -1.71%  '
$origStyle_D47 = $ws.Range('D47').Style
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.08'
$ws.Range('D47').Style = $origStyle_D47
$ws.Range('E47').Value = '  -3.94%  '
$origStyle_D48 = $ws.Range('D48').Style
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '87.71'
$ws.Range('D48').Style = $origStyle_D48
$ws.Range('E48').Value = '  +0.18%  '
$ws.Range('D49').Value = '0.0₆0100'
$ws.Range('E49').Value = '  -2.79%  '
$origStyle_D50 = $ws.Range('D50').Style
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0969'
$ws.Range('D50').Style = $origStyle_D50
$ws.Range('E50').Value = '  -3.70%  '
$origStyle_D51 = $ws.Range('D51').Style
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0497'
$ws.Range('D51').Style = $origStyle_D51
$ws.Range('E51').Value = '  -0.88%  '
